# ---------------------------------------------------------------------------
# Edit script: updates WR_89719272_WeekEnding_062925.xlsx per commit diff
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# -----------------------------------------------------------------------
# 1. Simple header / summary value updates
# -----------------------------------------------------------------------
$ws.Range("D5").Value = "Report Generated On: 08/18/2025 09:47 PM"
$ws.Range("C8").Value = 0
$ws.Range("C9").Value = 15

# -----------------------------------------------------------------------
# 2. Thursday section (rows 14-18): add a second data row.
#    Old layout:  16=Point35/POL-40-2 (data), 17=TOTAL
#    New layout:  16=Point34/POL-40-4 (data, replaces old content),
#                 17=Point35/POL-40-2 (data, new row - same values as old 16
#                     but priced at 0 and styled as an alternate-stripe row),
#                 18=TOTAL
# -----------------------------------------------------------------------

# Insert a new row before the old TOTAL row (row 17) to make room.
$ws.Rows("17:17").Insert()

# Clone the "alternate stripe" visual format (gray fill) from the Friday
# section's second data row (old row 23 -> now row 24 after the insert above)
# onto the newly inserted Thursday row 17.
$ws.Range("A24:H24").Copy()
$ws.Range("A17:H17").PasteSpecial(-4122)   # xlPasteFormats

# Row 16 becomes "Point 34 / POL-40-4 / Rem" with zeroed quantities.
$ws.Range("A16").Value = "Point 34"
$ws.Range("B16").Value = "POL-40-4"
$ws.Range("C16").Value = "Rem"
$ws.Range("D16").Value = "Pole,40ft,Class 4"
$ws.Range("E16").Value = "EA"
$ws.Range("F16").Value = 0
$ws.Range("H16").Value = 0

# Row 17 (new) carries what used to be row 16's content, now priced at 0.
$ws.Range("A17").Value = "Point 35"
$ws.Range("B17").Value = "POL-40-2"
$ws.Range("C17").Value = "Inst"
$ws.Range("D17").Value = "Pole,40ft,Class 2"
$ws.Range("E17").Value = "EA"
$ws.Range("F17").Value = 1
$ws.Range("H17").Value = 0

# TOTAL row (now row 18) pricing -> 0
$ws.Range("H18").Value = 0

# -----------------------------------------------------------------------
# 3. Friday section totals / pricing -> 0 (rows shifted down by 1 due to
#    the insert above: old 20-26 now 21-27)
# -----------------------------------------------------------------------
$ws.Range("H23").Value = 0
$ws.Range("H24").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("H26").Value = 0
$ws.Range("H27").Value = 0

# -----------------------------------------------------------------------
# 4. Saturday section (old rows 29-36, now 30-37 after the Thursday insert):
#    pricing -> 0 for existing rows, plus 4 brand new rows added before the
#    TOTAL row for "Point 35/POL-40-2" (qty 0) and three "Point 37" lines.
# -----------------------------------------------------------------------
$ws.Range("H32").Value = 0
$ws.Range("H33").Value = 0
$ws.Range("H34").Value = 0
$ws.Range("H35").Value = 0
$ws.Range("H36").Value = 0

# Insert 4 new rows before the Saturday TOTAL row (currently row 37).
$ws.Rows("37:40").Insert()

# Clone visual formats for the 4 new rows from the existing alternating
# stripe pattern directly above them (rows 34 "white" / 35 "gray" give us
# both flavors to copy from).
$ws.Range("A35:H35").Copy()
$ws.Range("A37:H37").PasteSpecial(-4122)   # xlPasteFormats (gray stripe)

$ws.Range("A34:H34").Copy()
$ws.Range("A38:H38").PasteSpecial(-4122)   # xlPasteFormats (white stripe)

$ws.Range("A35:H35").Copy()
$ws.Range("A39:H39").PasteSpecial(-4122)   # xlPasteFormats (gray stripe)

$ws.Range("A34:H34").Copy()
$ws.Range("A40:H40").PasteSpecial(-4122)   # xlPasteFormats (white stripe)

# Row 37: Point 35 / POL-40-2
$ws.Range("A37").Value = "Point 35"
$ws.Range("B37").Value = "POL-40-2"
$ws.Range("C37").Value = "Inst"
$ws.Range("D37").Value = "Pole,40ft,Class 2"
$ws.Range("E37").Value = "EA"
$ws.Range("F37").Value = 0
$ws.Range("H37").Value = 0

# Row 38: Point 37 / INS-15-P-S-C
$ws.Range("A38").Value = "Point 37"
$ws.Range("B38").Value = "INS-15-P-S-C"
$ws.Range("C38").Value = "Inst"
$ws.Range("D38").Value = "INS,15kV,Pin,Silicon Polymer,Corr"
$ws.Range("E38").Value = "EA"
$ws.Range("F38").Value = 0
$ws.Range("H38").Value = 0

# Row 39: Point 37 / PIN-15-PTP-C
$ws.Range("A39").Value = "Point 37"
$ws.Range("B39").Value = "PIN-15-PTP-C"
$ws.Range("C39").Value = "Inst"
$ws.Range("D39").Value = "Pin,15kV,Pole top,Corrosive"
$ws.Range("E39").Value = "EA"
$ws.Range("F39").Value = 0
$ws.Range("H39").Value = 0

# Row 40: Point 37 / POL-40-2
$ws.Range("A40").Value = "Point 37"
$ws.Range("B40").Value = "POL-40-2"
$ws.Range("C40").Value = "Inst"
$ws.Range("D40").Value = "Pole,40ft,Class 2"
$ws.Range("E40").Value = "EA"
$ws.Range("F40").Value = 0
$ws.Range("H40").Value = 0

# TOTAL row (now row 41) pricing -> 0
$ws.Range("H41").Value = 0
